# Update ERA table1 sheet with corrected fisheries data (rows 2-80).
# Writes Art (A), Procent (B), Kumulativ Procent (C), Fiske (D) for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Nordhavsräka"
$ws.Cells.Item(2, 2).Value = 97.8
$ws.Cells.Item(2, 3).Value = 97.8
$ws.Cells.Item(2, 4).Value = "Räkfiske rist;  3a4"

$ws.Cells.Item(3, 1).Value = "Nordhavsräka"
$ws.Cells.Item(3, 2).Value = 90.09999999999999
$ws.Cells.Item(3, 3).Value = 90.09999999999999
$ws.Cells.Item(3, 4).Value = "Räkfiske tunnel och rist;  3a4"

$ws.Cells.Item(4, 1).Value = "Torsk"
$ws.Cells.Item(4, 2).Value = 3.4
$ws.Cells.Item(4, 3).Value = 93.5
$ws.Cells.Item(4, 4).Value = "Räkfiske tunnel och rist;  3a4"

$ws.Cells.Item(5, 1).Value = "Gråsej"
$ws.Cells.Item(5, 2).Value = 2.8
$ws.Cells.Item(5, 3).Value = 96.3
$ws.Cells.Item(5, 4).Value = "Räkfiske tunnel och rist;  3a4"

$ws.Cells.Item(6, 1).Value = "Havskräfta"
$ws.Cells.Item(6, 2).Value = 99.59999999999999
$ws.Cells.Item(6, 3).Value = 99.59999999999999
$ws.Cells.Item(6, 4).Value = "Fiske med kräftburar;  3a"

$ws.Cells.Item(7, 1).Value = "Hummer"
$ws.Cells.Item(7, 2).Value = 44.4
$ws.Cells.Item(7, 3).Value = 44.4
$ws.Cells.Item(7, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(8, 1).Value = "Makrill"
$ws.Cells.Item(8, 2).Value = 17.7
$ws.Cells.Item(8, 3).Value = 62.1
$ws.Cells.Item(8, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(9, 1).Value = "Krabbtaska"
$ws.Cells.Item(9, 2).Value = 14
$ws.Cells.Item(9, 3).Value = 76.09999999999999
$ws.Cells.Item(9, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(10, 1).Value = "Berggylta"
$ws.Cells.Item(10, 2).Value = 6.3
$ws.Cells.Item(10, 3).Value = 82.40000000000001
$ws.Cells.Item(10, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(11, 1).Value = "Skärsnultra"
$ws.Cells.Item(11, 2).Value = 3.3
$ws.Cells.Item(11, 3).Value = 85.7
$ws.Cells.Item(11, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(12, 1).Value = "Stensnultra"
$ws.Cells.Item(12, 2).Value = 3.1
$ws.Cells.Item(12, 3).Value = 88.8
$ws.Cells.Item(12, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(13, 1).Value = "Torsk"
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = 90.8
$ws.Cells.Item(13, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(14, 1).Value = "Äkta tunga"
$ws.Cells.Item(14, 2).Value = 1.4
$ws.Cells.Item(14, 3).Value = 92.2
$ws.Cells.Item(14, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(15, 1).Value = "Sill / strömming"
$ws.Cells.Item(15, 2).Value = 1.3
$ws.Cells.Item(15, 3).Value = 93.5
$ws.Cells.Item(15, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(16, 1).Value = "Piggvar"
$ws.Cells.Item(16, 2).Value = 1
$ws.Cells.Item(16, 3).Value = 94.5
$ws.Cells.Item(16, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(17, 1).Value = "Bleka / lyrtorsk"
$ws.Cells.Item(17, 2).Value = 1
$ws.Cells.Item(17, 3).Value = 95.5
$ws.Cells.Item(17, 4).Value = "Fiske med passiva redskap;  3a"

$ws.Cells.Item(18, 1).Value = "Havskräfta"
$ws.Cells.Item(18, 2).Value = 93.59999999999999
$ws.Cells.Item(18, 3).Value = 93.59999999999999
$ws.Cells.Item(18, 4).Value = "Bottentrål havskräfta och fisk;  3a21"

$ws.Cells.Item(19, 1).Value = "Fjärsing"
$ws.Cells.Item(19, 2).Value = 1.9
$ws.Cells.Item(19, 3).Value = 95.5
$ws.Cells.Item(19, 4).Value = "Bottentrål havskräfta och fisk;  3a21"

$ws.Cells.Item(20, 1).Value = "Torsk"
$ws.Cells.Item(20, 2).Value = 29
$ws.Cells.Item(20, 3).Value = 29
$ws.Cells.Item(20, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(21, 1).Value = "Kolja"
$ws.Cells.Item(21, 2).Value = 17.7
$ws.Cells.Item(21, 3).Value = 46.7
$ws.Cells.Item(21, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(22, 1).Value = "Rödtunga"
$ws.Cells.Item(22, 2).Value = 10
$ws.Cells.Item(22, 3).Value = 56.8
$ws.Cells.Item(22, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(23, 1).Value = "Havskräfta"
$ws.Cells.Item(23, 2).Value = 9.699999999999999
$ws.Cells.Item(23, 3).Value = 66.40000000000001
$ws.Cells.Item(23, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(24, 1).Value = "Gråsej"
$ws.Cells.Item(24, 2).Value = 9.300000000000001
$ws.Cells.Item(24, 3).Value = 75.7
$ws.Cells.Item(24, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(25, 1).Value = "Marulk"
$ws.Cells.Item(25, 2).Value = 7.3
$ws.Cells.Item(25, 3).Value = 83.09999999999999
$ws.Cells.Item(25, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(26, 1).Value = "Kummel"
$ws.Cells.Item(26, 2).Value = 5.7
$ws.Cells.Item(26, 3).Value = 88.8
$ws.Cells.Item(26, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(27, 1).Value = "Rödspätta"
$ws.Cells.Item(27, 2).Value = 5.3
$ws.Cells.Item(27, 3).Value = 94.09999999999999
$ws.Cells.Item(27, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(28, 1).Value = "Bleka / lyrtorsk"
$ws.Cells.Item(28, 2).Value = 1.2
$ws.Cells.Item(28, 3).Value = 95.3
$ws.Cells.Item(28, 4).Value = "Bottentrål fisk;  3a20"

$ws.Cells.Item(29, 1).Value = "Havskräfta"
$ws.Cells.Item(29, 2).Value = 79.3
$ws.Cells.Item(29, 3).Value = 79.3
$ws.Cells.Item(29, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(30, 1).Value = "Torsk"
$ws.Cells.Item(30, 2).Value = 6.9
$ws.Cells.Item(30, 3).Value = 86.2
$ws.Cells.Item(30, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(31, 1).Value = "Marulk"
$ws.Cells.Item(31, 2).Value = 2.7
$ws.Cells.Item(31, 3).Value = 88.90000000000001
$ws.Cells.Item(31, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(32, 1).Value = "Rödtunga"
$ws.Cells.Item(32, 2).Value = 2.5
$ws.Cells.Item(32, 3).Value = 91.40000000000001
$ws.Cells.Item(32, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(33, 1).Value = "Kolja"
$ws.Cells.Item(33, 2).Value = 1.5
$ws.Cells.Item(33, 3).Value = 92.90000000000001
$ws.Cells.Item(33, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(34, 1).Value = "Gråsej"
$ws.Cells.Item(34, 2).Value = 1.5
$ws.Cells.Item(34, 3).Value = 94.40000000000001
$ws.Cells.Item(34, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(35, 1).Value = "Rödspätta"
$ws.Cells.Item(35, 2).Value = 1.2
$ws.Cells.Item(35, 3).Value = 95.59999999999999
$ws.Cells.Item(35, 4).Value = "Bottentrål havskräfta och fisk;  3a20"

$ws.Cells.Item(36, 1).Value = "Havskräfta"
$ws.Cells.Item(36, 2).Value = 98.59999999999999
$ws.Cells.Item(36, 3).Value = 98.59999999999999
$ws.Cells.Item(36, 4).Value = "Bottentrål havskräfta rist;  3a"

$ws.Cells.Item(37, 1).Value = "Torsk"
$ws.Cells.Item(37, 2).Value = 68
$ws.Cells.Item(37, 3).Value = 68
$ws.Cells.Item(37, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(38, 1).Value = "Sjurygg"
$ws.Cells.Item(38, 2).Value = 14.9
$ws.Cells.Item(38, 3).Value = 82.90000000000001
$ws.Cells.Item(38, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(39, 1).Value = "Rödspätta"
$ws.Cells.Item(39, 2).Value = 6.5
$ws.Cells.Item(39, 3).Value = 89.40000000000001
$ws.Cells.Item(39, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(40, 1).Value = "Äkta tunga"
$ws.Cells.Item(40, 2).Value = 3.3
$ws.Cells.Item(40, 3).Value = 92.7
$ws.Cells.Item(40, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(41, 1).Value = "Skrubbskädda"
$ws.Cells.Item(41, 2).Value = 2.2
$ws.Cells.Item(41, 3).Value = 94.90000000000001
$ws.Cells.Item(41, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(42, 1).Value = "Slätvar"
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = 96.90000000000001
$ws.Cells.Item(42, 4).Value = "Passiva redskap (torsk); 22-24"

$ws.Cells.Item(43, 1).Value = "Ål"
$ws.Cells.Item(43, 2).Value = 20.8
$ws.Cells.Item(43, 3).Value = 20.8
$ws.Cells.Item(43, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(44, 1).Value = "Sill / strömming"
$ws.Cells.Item(44, 2).Value = 19.9
$ws.Cells.Item(44, 3).Value = 40.7
$ws.Cells.Item(44, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(45, 1).Value = "Lax"
$ws.Cells.Item(45, 2).Value = 18.7
$ws.Cells.Item(45, 3).Value = 59.4
$ws.Cells.Item(45, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(46, 1).Value = "Siklöja"
$ws.Cells.Item(46, 2).Value = 14.1
$ws.Cells.Item(46, 3).Value = 73.5
$ws.Cells.Item(46, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(47, 1).Value = "Sikar"
$ws.Cells.Item(47, 2).Value = 10.9
$ws.Cells.Item(47, 3).Value = 84.40000000000001
$ws.Cells.Item(47, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(48, 1).Value = "Abborre"
$ws.Cells.Item(48, 2).Value = 9.199999999999999
$ws.Cells.Item(48, 3).Value = 93.59999999999999
$ws.Cells.Item(48, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(49, 1).Value = "Torsk"
$ws.Cells.Item(49, 2).Value = 1.6
$ws.Cells.Item(49, 3).Value = 95.3
$ws.Cells.Item(49, 4).Value = "Fiske med övriga passiva redskap; 22-32"

$ws.Cells.Item(50, 1).Value = "Torsk"
$ws.Cells.Item(50, 2).Value = 92.90000000000001
$ws.Cells.Item(50, 3).Value = 92.90000000000001
$ws.Cells.Item(50, 4).Value = "Fiske med stormaskig bottentrål (torsk); 25-32"

$ws.Cells.Item(51, 1).Value = "Skrubbskädda"
$ws.Cells.Item(51, 2).Value = 5.7
$ws.Cells.Item(51, 3).Value = 98.59999999999999
$ws.Cells.Item(51, 4).Value = "Fiske med stormaskig bottentrål (torsk); 25-32"

$ws.Cells.Item(52, 1).Value = "Skrubbskädda"
$ws.Cells.Item(52, 2).Value = 31.6
$ws.Cells.Item(52, 3).Value = 31.6
$ws.Cells.Item(52, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(53, 1).Value = "Piggvar"
$ws.Cells.Item(53, 2).Value = 23.6
$ws.Cells.Item(53, 3).Value = 55.2
$ws.Cells.Item(53, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(54, 1).Value = "Torsk"
$ws.Cells.Item(54, 2).Value = 16.5
$ws.Cells.Item(54, 3).Value = 71.7
$ws.Cells.Item(54, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(55, 1).Value = "Abborre"
$ws.Cells.Item(55, 2).Value = 15.7
$ws.Cells.Item(55, 3).Value = 87.40000000000001
$ws.Cells.Item(55, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(56, 1).Value = "Gädda"
$ws.Cells.Item(56, 2).Value = 6.5
$ws.Cells.Item(56, 3).Value = 93.90000000000001
$ws.Cells.Item(56, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(57, 1).Value = "Sikar"
$ws.Cells.Item(57, 2).Value = 2.6
$ws.Cells.Item(57, 3).Value = 96.40000000000001
$ws.Cells.Item(57, 4).Value = "Passiva redskap (torsk); 25-32"

$ws.Cells.Item(58, 1).Value = "Siklöja"
$ws.Cells.Item(58, 2).Value = 97.2
$ws.Cells.Item(58, 3).Value = 97.2
$ws.Cells.Item(58, 4).Value = "Fiske med finmaskig bottentrål efter pelagiska arter; 30-31"

$ws.Cells.Item(59, 1).Value = "Sill / strömming"
$ws.Cells.Item(59, 2).Value = 84.3
$ws.Cells.Item(59, 3).Value = 84.3
$ws.Cells.Item(59, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 21-24"

$ws.Cells.Item(60, 1).Value = "Skarpsill"
$ws.Cells.Item(60, 2).Value = 14.2
$ws.Cells.Item(60, 3).Value = 98.5
$ws.Cells.Item(60, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 21-24"

$ws.Cells.Item(61, 1).Value = "Torsk"
$ws.Cells.Item(61, 2).Value = 32.7
$ws.Cells.Item(61, 3).Value = 32.7
$ws.Cells.Item(61, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(62, 1).Value = "Gråsej"
$ws.Cells.Item(62, 2).Value = 30
$ws.Cells.Item(62, 3).Value = 62.6
$ws.Cells.Item(62, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(63, 1).Value = "Marulk"
$ws.Cells.Item(63, 2).Value = 13.4
$ws.Cells.Item(63, 3).Value = 76
$ws.Cells.Item(63, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(64, 1).Value = "Bleka / lyrtorsk"
$ws.Cells.Item(64, 2).Value = 5.9
$ws.Cells.Item(64, 3).Value = 81.90000000000001
$ws.Cells.Item(64, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(65, 1).Value = "Kolja"
$ws.Cells.Item(65, 2).Value = 4.9
$ws.Cells.Item(65, 3).Value = 86.8
$ws.Cells.Item(65, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(66, 1).Value = "Kummel"
$ws.Cells.Item(66, 2).Value = 2.8
$ws.Cells.Item(66, 3).Value = 89.59999999999999
$ws.Cells.Item(66, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(67, 1).Value = "Långa"
$ws.Cells.Item(67, 2).Value = 2.4
$ws.Cells.Item(67, 3).Value = 92
$ws.Cells.Item(67, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(68, 1).Value = "Hälleflundra"
$ws.Cells.Item(68, 2).Value = 1.8
$ws.Cells.Item(68, 3).Value = 93.8
$ws.Cells.Item(68, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(69, 1).Value = "Havskatter"
$ws.Cells.Item(69, 2).Value = 1.7
$ws.Cells.Item(69, 3).Value = 95.5
$ws.Cells.Item(69, 4).Value = "Bottentrål fisk;  4"

$ws.Cells.Item(70, 1).Value = "Torsk"
$ws.Cells.Item(70, 2).Value = 93.2
$ws.Cells.Item(70, 3).Value = 93.2
$ws.Cells.Item(70, 4).Value = "Fiske med stormaskig bottentrål (torsk); 22-24"

$ws.Cells.Item(71, 1).Value = "Rödspätta"
$ws.Cells.Item(71, 2).Value = 3.6
$ws.Cells.Item(71, 3).Value = 96.8
$ws.Cells.Item(71, 4).Value = "Fiske med stormaskig bottentrål (torsk); 22-24"

$ws.Cells.Item(72, 1).Value = "Skarpsill"
$ws.Cells.Item(72, 2).Value = 58.1
$ws.Cells.Item(72, 3).Value = 58.1
$ws.Cells.Item(72, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 25-29"

$ws.Cells.Item(73, 1).Value = "Sill / strömming"
$ws.Cells.Item(73, 2).Value = 41.4
$ws.Cells.Item(73, 3).Value = 99.59999999999999
$ws.Cells.Item(73, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 25-29"

$ws.Cells.Item(74, 1).Value = "Sill / strömming"
$ws.Cells.Item(74, 2).Value = 70.8
$ws.Cells.Item(74, 3).Value = 70.8
$ws.Cells.Item(74, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 3a4"

$ws.Cells.Item(75, 1).Value = "Skarpsill"
$ws.Cells.Item(75, 2).Value = 14.5
$ws.Cells.Item(75, 3).Value = 85.3
$ws.Cells.Item(75, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 3a4"

$ws.Cells.Item(76, 1).Value = "Makrill"
$ws.Cells.Item(76, 2).Value = 9.800000000000001
$ws.Cells.Item(76, 3).Value = 95.09999999999999
$ws.Cells.Item(76, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 3a4"

$ws.Cells.Item(77, 1).Value = "Sill / strömming"
$ws.Cells.Item(77, 2).Value = 98.09999999999999
$ws.Cells.Item(77, 3).Value = 98.09999999999999
$ws.Cells.Item(77, 4).Value = "Pelagiskt fiske med aktiva redskap (flyttrål, vad); 30-31"

$ws.Cells.Item(78, 1).Value = "Tobisfiskar"
$ws.Cells.Item(78, 2).Value = 99.59999999999999
$ws.Cells.Item(78, 3).Value = 99.59999999999999
$ws.Cells.Item(78, 4).Value = "Pelagiskt fiske med aktiva redskap (bottentrål);  3a204"

$ws.Cells.Item(79, 1).Value = "Sill / strömming"
$ws.Cells.Item(79, 2).Value = 60.2
$ws.Cells.Item(79, 3).Value = 60.2
$ws.Cells.Item(79, 4).Value = "Fiske med finmaskig bottentrål efter pelagiska arter; 25-29"

$ws.Cells.Item(80, 1).Value = "Skarpsill"
$ws.Cells.Item(80, 2).Value = 39.3
$ws.Cells.Item(80, 3).Value = 99.40000000000001
$ws.Cells.Item(80, 4).Value = "Fiske med finmaskig bottentrål efter pelagiska arter; 25-29"
